$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.250.95"
$ws.Range("E2").Value = "  +2.95%  "

$ws.Range("D3").Value = "1.718.92"
$ws.Range("E3").Value = "  +3.39%  "

$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'239.77"
$ws.Range("E5").Value = "  +1.20%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "'0.4725"
$ws.Range("E7").Value = "  -1.64%  "

$ws.Range("D8").Value = "'0.2623"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "'0.06202"
$ws.Range("E9").Value = "  +0.76%  "

$ws.Range("D10").Value = "1.715.76"
$ws.Range("E10").Value = "  +3.23%  "

$ws.Range("D11").Value = "'0.07070"
$ws.Range("E11").Value = "  -0.43%  "

$ws.Range("D12").Value = "'15.33"
$ws.Range("E12").Value = "  +3.72%  "

$ws.Range("D13").Value = "'0.5975"
$ws.Range("E13").Value = "  +1.85%  "

$ws.Range("D14").Value = "'4.424"
$ws.Range("E14").Value = "  +1.38%  "

$ws.Range("D15").Value = "'76.19"
$ws.Range("E15").Value = "  +2.03%  "

$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "'0.9998"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "26.256.97"
$ws.Range("E18").Value = "  +3.04%  "

$ws.Range("D19").Value = "'0.000006812"
$ws.Range("E19").Value = "  +1.14%  "

$ws.Range("D20").Value = "'11.54"

$ws.Range("D21").Value = "1.934.22"
$ws.Range("E21").Value = "  +3.34%  "

$ws.Range("D22").Value = "'4.545"
$ws.Range("E22").Value = "  +2.66%  "

$ws.Range("D23").Value = "'8.725"
$ws.Range("E23").Value = "  +0.68%  "

$ws.Range("D24").Value = "'5.269"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").Value = "'135.02"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("D26").Value = "'15.19"
$ws.Range("E26").Value = "  +1.12%  "

$ws.Range("D27").Value = "'1.400"
$ws.Range("E27").Value = "  +1.40%  "

$ws.Range("D28").Value = "'1.761"
$ws.Range("E28").Value = "  +3.03%  "

$ws.Range("D29").Value = "'107.42"
$ws.Range("E29").Value = "  +2.35%  "

$ws.Range("D30").Value = "'3.973"
$ws.Range("E30").Value = "  +0.83%  "

$ws.Range("D31").Value = "'3.676"
$ws.Range("E31").Value = "  +0.38%  "

$ws.Range("D32").Value = "'0.07757"
$ws.Range("E32").Value = "  +1.30%  "

$ws.Range("D33").Value = "'0.04441"
$ws.Range("E33").Value = "  +5.30%  "

$ws.Range("E34").Value = "  +0.26%  "

$ws.Range("D35").Value = "'0.9750"
$ws.Range("E35").Value = "  +2.62%  "

$ws.Range("D36").Value = "'0.6177"
$ws.Range("E36").Value = "  +1.32%  "

$ws.Range("D37").Value = "'0.9298"
$ws.Range("E37").Value = "  +7.37%  "

$ws.Range("D38").Value = "'111.72"
$ws.Range("E38").Value = "  +15.46%  "

$ws.Range("D39").Value = "'2.407"
$ws.Range("E39").Value = "  -7.40%  "

$ws.Range("D40").Value = "'1.917"
$ws.Range("E40").Value = "  +3.52%  "

$ws.Range("D41").Value = "'0.9999"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").Value = "'0.01477"
$ws.Range("E42").Value = "  +0.67%  "

$ws.Range("D43").Value = "'5.418"
$ws.Range("E43").Value = "  +14.28%  "

$ws.Range("D44").Value = "'0.3823"
$ws.Range("E44").Value = "  +1.85%  "

$ws.Range("D45").Value = "'0.1179"
$ws.Range("E45").Value = "  +4.91%  "

$ws.Range("D46").Value = "'6.262"
$ws.Range("E46").Value = "  +0.98%  "

$ws.Range("D47").Value = "'0.05260"
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'30.19"
$ws.Range("E48").Value = "  +1.88%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.709"
$ws.Range("E49").Value = "  +5.89%  "

$ws.Range("E50").Value = "  +1.76%  "

$ws.Range("E51").Value = "  +1.54%  "
